$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header numbering): add new header cells L1, M1 continuing the 0-based sequence ---
$ws.Range("L1").Value = 10
$ws.Range("M1").Value = 11
# Match the styling of the existing header cells (e.g. K1)
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2: Astronomy 201 ---
# Day-of-week code split out of the course title into its own column
$ws.Range("G2").Value = "TR"
$ws.Range("H2").Value = "Intro to Astrophysics: GW"
# Time range split into start / end time columns (force text so "0930" keeps its leading zero)
$ws.Range("I2").Formula = "=""1100"""
$ws.Range("I2").Copy()
$ws.Range("I2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("J2").Value = "1215pm"
# Shift building/room/instructor over by two columns
$ws.Range("K2").Value = "SCTR"
$ws.Range("L2").Value = "W209"
$ws.Range("M2").Value = "Stinebring Daniel"

# --- Row 3: Astronomy 302 ---
$ws.Range("G3").Value = "TR"
$ws.Range("H3").Value = "Astrophysics II"
$ws.Range("I3").Formula = "=""0930"""
$ws.Range("I3").Copy()
$ws.Range("I3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("J3").Value = "1050am"
$ws.Range("K3").Value = "SCTR"
$ws.Range("L3").Value = "W209"
$ws.Range("M3").Value = "Scudder Jillian"
